$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear E1 value while keeping its style
$ws.Cells.Item(1, 5).Value = $null

# Row 2
$ws.Range("B2").Value = "(`"electroencephalogram`" OR `"EEG`") AND (`"Bispectral index`" OR `"BIS`")"
$ws.Range("D2").Value = "01 October 2024 to 30 September 2025"
$ws.Range("F2").Value = "31.0, 31"

# Row 3
$ws.Range("B3").Value = "`"patient monitoring`" AND (`"anesthesia`" OR `"sedation`")"
$ws.Range("D3").Value = "01 October 2024 to 30 September 2025"
$ws.Range("F3").Value = "23.0, 23"

# Row 4
$ws.Range("B4").Value = "(`"anesthesia`" OR `"sedation`") AND `"EEG`""
$ws.Range("D4").Value = "01 October 2024 to 30 September 2025"
$ws.Range("F4").Value = "129.0, 129"

# Row 5
$ws.Range("B5").Value = "(`"EEG`" OR `"EEG module`" OR `"electroencephalogram`") AND `"Patient Monitoring`""
$ws.Range("D5").Value = "01 October 2024 to 30 September 2025"
$ws.Range("F5").Value = "4.0, 4"

# Row 6
$ws.Range("B6").Value = "`"Bispectral index`" AND `"hospital`" AND `"monitors`""
$ws.Range("D6").Value = "01 October 2024 to 30 September 2025"
$ws.Range("F6").Value = "19.0, 19"

# Row 7
$ws.Range("B7").Value = "(`"BIS`" OR `"Bispectral Index`") AND `"Philips`""
$ws.Range("D7").Value = "01 October 2024 to 30 September 2025"
$ws.Range("F7").Value = "2.0, 2"

# Row 8
$ws.Range("B8").Value = "(`"E BIS`" OR `"Bispectral Index`") AND `"GE Healthcare`""
$ws.Range("D8").Value = "01 October 2024 to 30 September 2025"
$ws.Range("F8").Value = "1.0, 1"

# Row 9
$ws.Range("B9").Value = "`"6800-30-50486`" AND `"Mindray`""
$ws.Range("D9").Value = "01 October 2024 to 30 September 2025"
$ws.Range("F9").Value = "0.0, 0"

# Row 10
$ws.Range("B10").Value = "`"Spacelabs`" AND `"Bispectral index`""
$ws.Range("F10").Value = "'0"
$ws.Range("F10").Style = "Normal"

# Row 11
$ws.Range("B11").Value = "`"Spacelabs`" AND `"91482`""
$ws.Range("F11").Value = "'0"
$ws.Range("F11").Style = "Normal"

# Row 12
$ws.Range("B12").Value = "(`"Bispectral Index`" OR `"Bispectral Index Analysis Module`") AND (`"brain data acquisition`" OR `"EEG data acquisition`" OR `"acquisition of EEG signals`" OR `"incidence of awareness`" OR `"anesthetic administration`" OR `"general anesthesia`" OR `"sedation`")"
$ws.Range("F12").Value = "'71"
$ws.Range("F12").Style = "Normal"

# Row 13
$ws.Range("B13").Value = "(`"Bispectral Index`") AND (`"brain`" OR `"data acquisition`" OR `"EEG signals`" OR `"anesthetic agent`")"
$ws.Range("F13").Value = "'16"
$ws.Range("F13").Style = "Normal"

# Row 14
$ws.Range("B14").Value = "(`"Bispectral Index`") AND (`"no display`" OR `"loss of data`" OR `"no alarm`" OR `"injury`") AND `"risk`""
$ws.Range("F14").Value = "'0"
$ws.Range("F14").Style = "Normal"

# Row 15
$ws.Range("B15").Value = "(`"Bispectral Index`") AND (`"BISx Interface Cable`" OR `"BISx`""
$ws.Range("F15").Value = "'0"
$ws.Range("F15").Style = "Normal"

# Row 16
$ws.Range("B16").Value = "`"BIS pod`" OR `"BIS Sensor`""
$ws.Range("F16").Value = "'5"
$ws.Range("F16").Style = "Normal"

# Row 17
$ws.Range("B17").Value = "(`"Bispectral Index`") AND (`"real-time encephalogram`" OR `"real-time EEG`") AND (`"adults`" AND `"pediatrics`")"
$ws.Range("F17").Value = "'0"
$ws.Range("F17").Style = "Normal"

# Row 18
$ws.Range("B18").Value = "(`"brain wave measurement`" OR `"brain wave frequency`" OR `"depth of consciousness`") AND (`"EEG`" OR `"BIS`")"
$ws.Range("F18").Value = "'2"
$ws.Range("F18").Style = "Normal"

# Row 19
$ws.Range("B19").Value = "(`"supression ratio`" OR `"Spectral Edge Frequency`" OR `"Median Power Frequency`" OR `"Electromyographic strength`" OR `"Signal Quality Index`" OR `"Burst Count`") AND `"BIS`""
$ws.Range("F19").Value = "'2"
$ws.Range("F19").Style = "Normal"

# Row 20
$ws.Range("B20").Value = "(`"depth of consciousness`" OR `"sedation level`") AND `"BIS`""
$ws.Range("F20").Value = "'1"
$ws.Range("F20").Style = "Normal"

# Column widths (closest achievable values given the engine's internal
# pixel-based rounding of ColumnWidth; targets are 132.42578125 and 39.7109375)
$ws.Columns.Item(3).ColumnWidth = 131.67
$ws.Columns.Item(5).ColumnWidth = 38.835

# Selection
$ws.Range("E15").Select()

# Page margins (Excel defaults in inches: 0.75/0.75/1/1, header/footer 0.5/0.5)
$ws.PageSetup.LeftMargin = $excel.InchesToPoints(0.75)
$ws.PageSetup.RightMargin = $excel.InchesToPoints(0.75)
$ws.PageSetup.TopMargin = $excel.InchesToPoints(1)
$ws.PageSetup.BottomMargin = $excel.InchesToPoints(1)
$ws.PageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
$ws.PageSetup.FooterMargin = $excel.InchesToPoints(0.5)
